$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17, shifting the existing data (rows 17-137)
# down to rows 18-138.
$ws.Rows("17").Insert()

# Populate the newly inserted row 17 with the new weekly price-report entry.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R are constant across this market/category
# subset, so they are copied from the surrounding rows; D,J,K,L,M,P carry
# the new observation's values.
$ws.Range("A17").Value = 3
$ws.Range("B17").Value = "Femacal de La Calera"
$ws.Range("C17").Value = "Coquimbo"
$ws.Range("D17").Value = 44602
$ws.Range("E17").Value = 5
$ws.Range("F17").Value = 100112052
$ws.Range("G17").Value = "Albahaca"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 130
$ws.Range("K17").Value = 4000
$ws.Range("L17").Value = 4500
$ws.Range("M17").Value = 4231
$ws.Range("N17").Value = "$/docena de matas"
$ws.Range("O17").Value = "Provincia de Quillota"
$ws.Range("P17").Value = 705
$ws.Range("Q17").Value = 6
$ws.Range("R17").Value = "Hortaliza"

# Match the date-format style used by the other rows' "Fecha" column (D).
$ws.Range("D17").NumberFormat = $ws.Range("D18").NumberFormat
